$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data; force text format on numeric-looking
# strings (e.g. "20.231.37", "0.06536") so Excel does not auto-convert them to numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "20.231.37"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.48%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.438.80"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.04%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.012"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.97%  "
$ws.Range("B5").Value = "USDC"
$ws.Range("C5").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9089"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -9.25%  "
$ws.Range("B6").Value = "BNB"
$ws.Range("C6").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "277.83"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.48%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.36%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3132"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.73%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "38.80"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.89%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.018"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.19%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06536"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.18%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.006"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.44%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.406"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.10%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "17.54"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +3.44%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.069"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.23%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.443.42"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.35%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001019"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.19%  "
$ws.Range("B18").Value = "Dai"
$ws.Range("C18").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9119"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -8.92%  "
$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.05612"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.32%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "67.33"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -8.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.410"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -3.38%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.41"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.27%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.78%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.251"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.20%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "20.307.71"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.82%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.200"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.14%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "135.24"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.61%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.89"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.80%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.595.03"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.66%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "110.54"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.47%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.718"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.12%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.8181"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.46%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.880"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -8.36%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.07621"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.84%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.05965"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +3.53%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.471"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +15.37%  "
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.151"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +8.47%  "
$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.703"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.09%  "
$ws.Range("B39").Value = "Aptos"
$ws.Range("C39").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "10.27"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.14%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01995"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.13%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1827"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -5.56%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9179"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -8.26%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.529"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5232"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.84%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.834"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -18.78%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "11.89"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.64%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "119.52"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +7.57%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5148"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.87%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.764"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.91%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06340"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.87%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9933"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.74%  "

# Restore the cells original (default) style now that the text values are set
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Style = "Normal"

